# novas traducoes e melhorias
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GlobalTile")

# New / improved PTBR translations
$ws.Range("C4").Value = "Planície"        # id 1 (Plain): Plano -> Planície
$ws.Range("C13").Value = "Penhasco Nevado" # id 9 (Cliff, snow variant): Penhasco -> Penhasco Nevado
$ws.Range("C15").Value = "Estrada"         # id 11 (road): estrada -> Estrada

# Column C widened to fit the longer text (best-fit)
$ws.Range("C:C").EntireColumn.AutoFit()

# Selection moved to C15
$ws.Range("C15").Select()
